$d = $word.ActiveDocument

# 1. Correct school spelling: "North Central University" -> "Northcentral University"
$d.Content.Find.Execute("North Central University", $true, $false, $false, $false, $false, $true, 1, $false, "Northcentral University", 2) | Out-Null

# 2. Word-level edits to paragraph about media/escalation perspectives
$d.Content.Find.Execute("While these perspectives effectively drive", $true, $false, $false, $false, $false, $true, 1, $false, "While these perspectives efficiently drive", 2) | Out-Null
$d.Content.Find.Execute("explain the need for forceful escalations", $true, $false, $false, $false, $false, $true, 1, $false, "explain the need for violent escalations", 2) | Out-Null
